$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 333.7143
$ws.Range("I8").Value = 139.33333
$ws.Range("J8").Value = 1500
$ws.Range("K8").Value = 417.99999
$ws.Range("L8").Value = 4500
$ws.Range("M8").Value = -278.99999
$ws.Range("N8").Value = -4778
$ws.Range("H46").Value = 1526.8
$ws.Range("I46").Value = 817
$ws.Range("J46").Value = 2000
$ws.Range("K46").Value = 2451
$ws.Range("L46").Value = 6000
$ws.Range("M46").Value = -2332
$ws.Range("N46").Value = -6238
$ws.Range("H60").Value = 1526.8
$ws.Range("I60").Value = 817
$ws.Range("J60").Value = 2000
$ws.Range("K60").Value = 2451
$ws.Range("L60").Value = 6000
$ws.Range("M60").Value = -1967
$ws.Range("N60").Value = -6968
$ws.Range("H86").Value = 2227.6667
$ws.Range("I86").Value = 2227.6667
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 2227.6667
$ws.Range("L86").Value = ""
$ws.Range("M86").Value = -1104.6667
$ws.Range("N86").Value = ""
$ws.Range("H89").Value = 2227.6667
$ws.Range("I89").Value = 2227.6667
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 11138.3335
$ws.Range("L89").Value = ""
$ws.Range("M89").Value = -5522.333500000001
$ws.Range("N89").Value = ""
$ws.Range("H100").Value = 40004396
$ws.Range("J100").Value = 7000
$ws.Range("L100").Value = 7000
$ws.Range("N100").Value = -8082
$ws.Range("H123").Value = 42302.855
$ws.Range("J123").Value = 42302.855
$ws.Range("L123").Value = 42302.855
$ws.Range("N123").Value = -52102.855
$ws.Range("H138").Value = 2778.5
$ws.Range("I138").Value = 1716.1875
$ws.Range("J138").Value = 3722.7778
$ws.Range("K138").Value = 5148.5625
$ws.Range("L138").Value = 11168.3334
$ws.Range("M138").Value = -8.5625
$ws.Range("N138").Value = -21448.3334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 399
$ws.Range("I5").Value = 300
$ws.Range("K5").Value = 300
$ws.Range("M5").Value = -188
$ws.Range("H97").Value = 1122.8572
$ws.Range("I97").Value = 982.7273
$ws.Range("K97").Value = 982.7273
$ws.Range("M97").Value = -486.7273
$ws.Range("H102").Value = 1941.8
$ws.Range("I102").Value = 1941.8
$ws.Range("K102").Value = 1941.8
$ws.Range("M102").Value = -319.8
$ws.Range("H132").Value = 2359.4075
$ws.Range("I132").Value = 1854.2609
$ws.Range("J132").Value = 5264
$ws.Range("K132").Value = 5562.7827
$ws.Range("L132").Value = 15792
$ws.Range("M132").Value = -3032.7827
$ws.Range("N132").Value = -20852

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 399
$ws.Range("I4").Value = 300
$ws.Range("K4").Value = 300
$ws.Range("M4").Value = -185
$ws.Range("H15").Value = 32000
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 32000
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = ""
$ws.Range("M15").Value = ""
$ws.Range("N15").Value = -32454
$ws.Range("H19").Value = 30009.5
$ws.Range("J19").Value = 30009.5
$ws.Range("L19").Value = 30009.5
$ws.Range("N19").Value = -30355.5
$ws.Range("H82").Value = 22027.5
$ws.Range("J82").Value = 32992
$ws.Range("L82").Value = 32992
$ws.Range("N82").Value = -33758
$ws.Range("H85").Value = 22027.5
$ws.Range("J85").Value = 32992
$ws.Range("L85").Value = 32992
$ws.Range("N85").Value = -35644
$ws.Range("H105").Value = 1707.5555
$ws.Range("I105").Value = 1672.0571
$ws.Range("J105").Value = 2950
$ws.Range("K105").Value = 1672.0571
$ws.Range("L105").Value = 2950
$ws.Range("M105").Value = 74.94290000000001
$ws.Range("N105").Value = -6444

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3285
$ws.Range("I62").Value = 2999.1667
$ws.Range("K62").Value = 2999.1667
$ws.Range("M62").Value = -2375.1667
$ws.Range("H65").Value = 3285
$ws.Range("I65").Value = 2999.1667
$ws.Range("K65").Value = 14995.8335
$ws.Range("M65").Value = -11875.8335
$ws.Range("H105").Value = 1951.3636
$ws.Range("I105").Value = 1922
$ws.Range("J105").Value = 2014.2858
$ws.Range("K105").Value = 1922
$ws.Range("L105").Value = 2014.2858
$ws.Range("M105").Value = -175
$ws.Range("N105").Value = -5508.2858
$ws.Range("H132").Value = 3958.037
$ws.Range("I132").Value = 3168.6191
$ws.Range("J132").Value = 6721
$ws.Range("K132").Value = 9505.8573
$ws.Range("L132").Value = 20163
$ws.Range("M132").Value = -6975.8573
$ws.Range("N132").Value = -25223

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 3276.923
$ws.Range("J68").Value = 5480.731
$ws.Range("L68").Value = 16442.193
$ws.Range("N68").Value = -18064.193
$ws.Range("H71").Value = 3276.923
$ws.Range("J71").Value = 5480.731
$ws.Range("L71").Value = 49326.579
$ws.Range("N71").Value = -57438.579
$ws.Range("H118").Value = 2514.85
$ws.Range("I118").Value = 613
$ws.Range("J118").Value = 4070.9092
$ws.Range("K118").Value = 1839
$ws.Range("L118").Value = 12212.7276
$ws.Range("M118").Value = -596
$ws.Range("N118").Value = -14698.7276

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2826.9
$ws.Range("I80").Value = 2697.1428
$ws.Range("J80").Value = 3129.6667
$ws.Range("K80").Value = 2697.1428
$ws.Range("L80").Value = 3129.6667
$ws.Range("M80").Value = -1699.1428
$ws.Range("N80").Value = -5125.6667
$ws.Range("H83").Value = 2826.9
$ws.Range("I83").Value = 2697.1428
$ws.Range("J83").Value = 3129.6667
$ws.Range("K83").Value = 13485.714
$ws.Range("L83").Value = 15648.3335
$ws.Range("M83").Value = -8493.714
$ws.Range("N83").Value = -25632.3335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4324.375
$ws.Range("I7").Value = 2565.8333
$ws.Range("J7").Value = 9600
$ws.Range("K7").Value = 2565.8333
$ws.Range("L7").Value = 9600
$ws.Range("M7").Value = -2453.8333
$ws.Range("N7").Value = -9824
$ws.Range("H126").Value = 4324.375
$ws.Range("I126").Value = 2565.8333
$ws.Range("J126").Value = 9600
$ws.Range("K126").Value = 7697.499899999999
$ws.Range("L126").Value = 28800
$ws.Range("M126").Value = -5227.499899999999
$ws.Range("N126").Value = -33740

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 667.3333
$ws.Range("I107").Value = 505.25
$ws.Range("J107").Value = 916.6923
$ws.Range("K107").Value = 1515.75
$ws.Range("L107").Value = 2750.0769
$ws.Range("M107").Value = 404.25
$ws.Range("N107").Value = -6590.0769
$ws.Range("H136").Value = 2521.0425
$ws.Range("I136").Value = 1034.6552
$ws.Range("J136").Value = 4915.778
$ws.Range("K136").Value = 3103.9656
$ws.Range("L136").Value = 14747.334
$ws.Range("M136").Value = -553.9655999999995
$ws.Range("N136").Value = -19847.334

Write-Output "done"